$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing bound rows (11 and 12) towards planned capacities
$ws.Range("I11").Value = 2000
$ws.Range("J11").Value = 3000

$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0

# Row 13 - new CAP_BND LO entry (planned capacity bound, up)
$ws.Range("D13").Value = 2045
$ws.Range("E13").Value = "LO"
$ws.Range("F13").Value = "CAP_BND"
$ws.Range("I13").Value = 3000
$ws.Range("J13").Value = 3500
$ws.Range("M13").Value = "ERWINELCWIN3N"

# Row 14 - new CAP_BND LO entry (planned capacity bound, down)
$ws.Range("D14").Value = 2045
$ws.Range("E14").Value = "LO"
$ws.Range("F14").Value = "CAP_BND"
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 3500
$ws.Range("M14").Value = "ERWINELCWIN5N"

# Row 15 (previously row 13 shifted down) - C15 loses its Pset_PN text value
$ws.Range("C15").Value = $null

# Update selection to match new active cell
$ws.Range("F22").Select()
